$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI")

# Header row updates
$ws.Range("B1").Value = "What_Action_Items"
$ws.Range("D1").Value = "Who_Responsible"
$ws.Range("E1").Value = "When_Due_date"

# VoC project action items (row 8) - add deadlines
$ws.Range("F8").Value = "1. Need discussion with Yamazaki-san and other party to solve user's question by end of December`n2. Send one VoC survey after providing answers to users by end of January 2022"

# Knowledge base chatbot project action items (row 12) - add deadline
$ws.Range("F12").Value = "1. Have to discuss further with IT representative by mid of January`n"

# Update selected cell
$ws.Range("E1").Select()
